$wb = $excel.ActiveWorkbook

# Updated "想去人数" (column F) values for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    3  = 570
    6  = 97
    7  = 276
    8  = 72
    9  = 1180
    10 = 16675
    11 = 293
    12 = 211
    14 = 6456
    15 = 650
    16 = 134
    18 = 33
    21 = 68
    24 = 41
    25 = 25
    26 = 11
    28 = 236
    29 = 912
    30 = 69
    31 = 5076
    33 = 11469
    34 = 1253
    36 = 165
    37 = 225
    38 = 3861
}
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Updated "想去人数" (column F) values for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    3  = 570
    6  = 97
    7  = 276
    8  = 72
    9  = 1180
    10 = 16675
    11 = 293
    12 = 211
    14 = 6456
    15 = 650
    16 = 134
    18 = 33
    21 = 68
    24 = 41
    25 = 25
    26 = 11
    28 = 236
    29 = 912
    30 = 69
    31 = 5076
    34 = 11469
    35 = 1253
    37 = 165
    38 = 225
    39 = 3861
}
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
